$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("connected component")
$ws2 = $wb.Worksheets.Item("degree")

$ws1.Range("A9").Value = "the average degree"
$ws1.Range("B9").Value = 4.19753086419753

$ws1.Range("E14").Select() | Out-Null
$ws2.Range("E10").Select() | Out-Null
$ws1.Activate() | Out-Null
